# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Updates DAMSLTag (col I) and DialogAct (col J)
# for a set of rows whose sentences were re-classified.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 28;  DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 59;  DAMSLTag = "aa"; DialogAct = "Agree/Accept" },
    @{ Row = 65;  DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 91;  DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 104; DAMSLTag = "ba"; DialogAct = "Appreciation" },
    @{ Row = 107; DAMSLTag = "ba"; DialogAct = "Appreciation" },
    @{ Row = 109; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 112; DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 118; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 123; DAMSLTag = "sd"; DialogAct = "Statement-non-opinion" },
    @{ Row = 141; DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 146; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 158; DAMSLTag = "b";  DialogAct = "Acknowledge (Backchannel)" },
    @{ Row = 160; DAMSLTag = "sv"; DialogAct = "Statement-opinion" },
    @{ Row = 166; DAMSLTag = "aa"; DialogAct = "Agree/Accept" }
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 9).Value = $change.DAMSLTag
    $ws.Cells.Item($change.Row, 10).Value = $change.DialogAct
}
